$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 12 new data rows right after the table header (row 15), pushing
#    the existing "JESUS MANUEL" worker rows (old 16:20) and the footer rows
#    (old 25:26) down to make room for two new workers (LUIS JAIRO and
#    JUAN CARLOS), 6 period-rows each.
# ---------------------------------------------------------------------------
$ws.Rows("16:27").Insert()

# ---------------------------------------------------------------------------
# 2. The freshly inserted rows are blank/unformatted. Clone the formatting
#    (styles/borders/number formats) of the existing body-row pattern
#    (row 28, which is the former row 16 - a normal "middle" table row) onto
#    all 12 new rows.
# ---------------------------------------------------------------------------
$ws.Range("B28:J28").Copy()
$ws.Range("B16:J27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Fill in the worker data. Two new workers, 6 periods each.
# ---------------------------------------------------------------------------
$newWorkers = @(
    @("1049926522", "LUIS JAIRO BELTRAN SANJUAN", @(
        @("2507", 56940, 1423500),
        @("2506", 56940, 1423500),
        @("2505", 56940, 1423500),
        @("2504", 52000, 1423500),
        @("2502", 52000, 1423500),
        @("2501", 52000, 1423500)
    )),
    @("1049932844", "JUAN CARLOS NAVARRO SANCHEZ", @(
        @("2507", 56940, 1423500),
        @("2506", 56940, 1423500),
        @("2505", 56940, 1423500),
        @("2504", 52000, 1423500),
        @("2502", 52000, 1423500),
        @("2501", 52000, 1423500)
    ))
)

$r = 16
foreach ($worker in $newWorkers) {
    $docId = $worker[0]
    $name = $worker[1]
    $periods = $worker[2]
    foreach ($p in $periods) {
        $ws.Cells.Item($r, 2).Value = "CC"
        $ws.Cells.Item($r, 3).Value = $docId
        $ws.Cells.Item($r, 4).Value = $name
        $ws.Cells.Item($r, 5).Value = $p[0]
        $ws.Cells.Item($r, 6).Value = $p[1]
        $ws.Cells.Item($r, 7).Value = $p[2]
        $r++
    }
}

# ---------------------------------------------------------------------------
# 4. The pre-existing worker (JESUS MANUEL, now at rows 28:32) keeps the same
#    doc/name/amounts but the period order is reversed (most recent first).
# ---------------------------------------------------------------------------
$existingPeriods = @("1909", "1908", "1907", "1906", "1905")
$r = 28
foreach ($period in $existingPeriods) {
    $ws.Cells.Item($r, 5).Value = $period
    $r++
}

# ---------------------------------------------------------------------------
# 5. Update the summary header fields.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 819265
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 11
